$wb = $excel.ActiveWorkbook

$wsPIR = $wb.Worksheets.Item("PIR")
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsTemperature = $wb.Worksheets.Item("Temperature")
$wsProximity = $wb.Worksheets.Item("Proximity")

# --- PIR sheet: append rows 104-116 ---
$wsPIR.Cells.Item(104, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(104, 2).Value = "16:45:42"
$wsPIR.Cells.Item(104, 3).Value = "16:00"
$wsPIR.Cells.Item(104, 4).Value = "Bathroom"
$wsPIR.Cells.Item(104, 5).Value = "No Motion"
$wsPIR.Cells.Item(104, 6).Value = "Inactive"
$wsPIR.Cells.Item(105, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(105, 2).Value = "16:45:43"
$wsPIR.Cells.Item(105, 3).Value = "16:00"
$wsPIR.Cells.Item(105, 4).Value = "Bathroom"
$wsPIR.Cells.Item(105, 5).Value = "No Motion"
$wsPIR.Cells.Item(105, 6).Value = "Inactive"
$wsPIR.Cells.Item(106, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(106, 2).Value = "16:45:48"
$wsPIR.Cells.Item(106, 3).Value = "16:00"
$wsPIR.Cells.Item(106, 4).Value = "Bathroom"
$wsPIR.Cells.Item(106, 5).Value = "No Motion"
$wsPIR.Cells.Item(106, 6).Value = "Inactive"
$wsPIR.Cells.Item(107, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(107, 2).Value = "16:45:53"
$wsPIR.Cells.Item(107, 3).Value = "16:00"
$wsPIR.Cells.Item(107, 4).Value = "Bathroom"
$wsPIR.Cells.Item(107, 5).Value = "No Motion"
$wsPIR.Cells.Item(107, 6).Value = "Inactive"
$wsPIR.Cells.Item(108, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(108, 2).Value = "16:45:58"
$wsPIR.Cells.Item(108, 3).Value = "16:00"
$wsPIR.Cells.Item(108, 4).Value = "Bathroom"
$wsPIR.Cells.Item(108, 5).Value = "No Motion"
$wsPIR.Cells.Item(108, 6).Value = "Inactive"
$wsPIR.Cells.Item(109, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(109, 2).Value = "16:45:59"
$wsPIR.Cells.Item(109, 3).Value = "16:00"
$wsPIR.Cells.Item(109, 4).Value = "Bathroom"
$wsPIR.Cells.Item(109, 5).Value = "Motion Detected"
$wsPIR.Cells.Item(109, 6).Value = "Active"
$wsPIR.Cells.Item(110, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(110, 2).Value = "16:46:07"
$wsPIR.Cells.Item(110, 3).Value = "16:00"
$wsPIR.Cells.Item(110, 4).Value = "Bathroom"
$wsPIR.Cells.Item(110, 5).Value = "No Motion"
$wsPIR.Cells.Item(110, 6).Value = "Inactive"
$wsPIR.Cells.Item(111, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(111, 2).Value = "16:46:09"
$wsPIR.Cells.Item(111, 3).Value = "16:00"
$wsPIR.Cells.Item(111, 4).Value = "Bathroom"
$wsPIR.Cells.Item(111, 5).Value = "Motion Detected"
$wsPIR.Cells.Item(111, 6).Value = "Active"
$wsPIR.Cells.Item(112, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(112, 2).Value = "16:46:16"
$wsPIR.Cells.Item(112, 3).Value = "16:00"
$wsPIR.Cells.Item(112, 4).Value = "Bathroom"
$wsPIR.Cells.Item(112, 5).Value = "No Motion"
$wsPIR.Cells.Item(112, 6).Value = "Inactive"
$wsPIR.Cells.Item(113, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(113, 2).Value = "16:46:21"
$wsPIR.Cells.Item(113, 3).Value = "16:00"
$wsPIR.Cells.Item(113, 4).Value = "Bathroom"
$wsPIR.Cells.Item(113, 5).Value = "No Motion"
$wsPIR.Cells.Item(113, 6).Value = "Inactive"
$wsPIR.Cells.Item(114, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(114, 2).Value = "16:46:26"
$wsPIR.Cells.Item(114, 3).Value = "16:00"
$wsPIR.Cells.Item(114, 4).Value = "Bathroom"
$wsPIR.Cells.Item(114, 5).Value = "No Motion"
$wsPIR.Cells.Item(114, 6).Value = "Inactive"
$wsPIR.Cells.Item(115, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(115, 2).Value = "16:46:31"
$wsPIR.Cells.Item(115, 3).Value = "16:00"
$wsPIR.Cells.Item(115, 4).Value = "Bathroom"
$wsPIR.Cells.Item(115, 5).Value = "No Motion"
$wsPIR.Cells.Item(115, 6).Value = "Inactive"
$wsPIR.Cells.Item(116, 1).Value = "'2026-01-28"
$wsPIR.Cells.Item(116, 2).Value = "16:46:36"
$wsPIR.Cells.Item(116, 3).Value = "16:00"
$wsPIR.Cells.Item(116, 4).Value = "Bathroom"
$wsPIR.Cells.Item(116, 5).Value = "No Motion"
$wsPIR.Cells.Item(116, 6).Value = "Inactive"

# --- Humidity sheet: append rows 104-116 ---
$wsHumidity.Cells.Item(104, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(104, 2).Value = "16:45:41"
$wsHumidity.Cells.Item(104, 3).Value = "16:00"
$wsHumidity.Cells.Item(104, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(104, 5).Value = "'87.6%"
$wsHumidity.Cells.Item(104, 6).Value = "Active"
$wsHumidity.Cells.Item(105, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(105, 2).Value = "16:45:42"
$wsHumidity.Cells.Item(105, 3).Value = "16:00"
$wsHumidity.Cells.Item(105, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(105, 5).Value = "'87.6%"
$wsHumidity.Cells.Item(105, 6).Value = "Active"
$wsHumidity.Cells.Item(106, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(106, 2).Value = "16:45:44"
$wsHumidity.Cells.Item(106, 3).Value = "16:00"
$wsHumidity.Cells.Item(106, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(106, 5).Value = "'86.8%"
$wsHumidity.Cells.Item(106, 6).Value = "Active"
$wsHumidity.Cells.Item(107, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(107, 2).Value = "16:45:48"
$wsHumidity.Cells.Item(107, 3).Value = "16:00"
$wsHumidity.Cells.Item(107, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(107, 5).Value = "'87.7%"
$wsHumidity.Cells.Item(107, 6).Value = "Active"
$wsHumidity.Cells.Item(108, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(108, 2).Value = "16:45:56"
$wsHumidity.Cells.Item(108, 3).Value = "16:00"
$wsHumidity.Cells.Item(108, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(108, 5).Value = "'86.8%"
$wsHumidity.Cells.Item(108, 6).Value = "Active"
$wsHumidity.Cells.Item(109, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(109, 2).Value = "16:46:09"
$wsHumidity.Cells.Item(109, 3).Value = "16:00"
$wsHumidity.Cells.Item(109, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(109, 5).Value = "'88.0%"
$wsHumidity.Cells.Item(109, 6).Value = "Active"
$wsHumidity.Cells.Item(110, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(110, 2).Value = "16:46:17"
$wsHumidity.Cells.Item(110, 3).Value = "16:00"
$wsHumidity.Cells.Item(110, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(110, 5).Value = "'87.9%"
$wsHumidity.Cells.Item(110, 6).Value = "Active"
$wsHumidity.Cells.Item(111, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(111, 2).Value = "16:46:21"
$wsHumidity.Cells.Item(111, 3).Value = "16:00"
$wsHumidity.Cells.Item(111, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(111, 5).Value = "'87.8%"
$wsHumidity.Cells.Item(111, 6).Value = "Active"
$wsHumidity.Cells.Item(112, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(112, 2).Value = "16:46:25"
$wsHumidity.Cells.Item(112, 3).Value = "16:00"
$wsHumidity.Cells.Item(112, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(112, 5).Value = "'86.9%"
$wsHumidity.Cells.Item(112, 6).Value = "Active"
$wsHumidity.Cells.Item(113, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(113, 2).Value = "16:46:29"
$wsHumidity.Cells.Item(113, 3).Value = "16:00"
$wsHumidity.Cells.Item(113, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(113, 5).Value = "'87.8%"
$wsHumidity.Cells.Item(113, 6).Value = "Active"
$wsHumidity.Cells.Item(114, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(114, 2).Value = "16:46:33"
$wsHumidity.Cells.Item(114, 3).Value = "16:00"
$wsHumidity.Cells.Item(114, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(114, 5).Value = "'87.8%"
$wsHumidity.Cells.Item(114, 6).Value = "Active"
$wsHumidity.Cells.Item(115, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(115, 2).Value = "16:46:37"
$wsHumidity.Cells.Item(115, 3).Value = "16:00"
$wsHumidity.Cells.Item(115, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(115, 5).Value = "'86.8%"
$wsHumidity.Cells.Item(115, 6).Value = "Active"
$wsHumidity.Cells.Item(116, 1).Value = "'2026-01-28"
$wsHumidity.Cells.Item(116, 2).Value = "16:46:41"
$wsHumidity.Cells.Item(116, 3).Value = "16:00"
$wsHumidity.Cells.Item(116, 4).Value = "Bathroom"
$wsHumidity.Cells.Item(116, 5).Value = "'87.6%"
$wsHumidity.Cells.Item(116, 6).Value = "Active"

# --- Temperature sheet: append rows 104-115 ---
$wsTemperature.Cells.Item(104, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(104, 2).Value = "16:45:41"
$wsTemperature.Cells.Item(104, 3).Value = "16:00"
$wsTemperature.Cells.Item(104, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(104, 5).Value = "22.9C"
$wsTemperature.Cells.Item(104, 6).Value = "Active"
$wsTemperature.Cells.Item(105, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(105, 2).Value = "16:45:42"
$wsTemperature.Cells.Item(105, 3).Value = "16:00"
$wsTemperature.Cells.Item(105, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(105, 5).Value = "22.9C"
$wsTemperature.Cells.Item(105, 6).Value = "Active"
$wsTemperature.Cells.Item(106, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(106, 2).Value = "16:45:45"
$wsTemperature.Cells.Item(106, 3).Value = "16:00"
$wsTemperature.Cells.Item(106, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(106, 5).Value = "22.9C"
$wsTemperature.Cells.Item(106, 6).Value = "Active"
$wsTemperature.Cells.Item(107, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(107, 2).Value = "16:45:49"
$wsTemperature.Cells.Item(107, 3).Value = "16:00"
$wsTemperature.Cells.Item(107, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(107, 5).Value = "22.9C"
$wsTemperature.Cells.Item(107, 6).Value = "Active"
$wsTemperature.Cells.Item(108, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(108, 2).Value = "16:45:57"
$wsTemperature.Cells.Item(108, 3).Value = "16:00"
$wsTemperature.Cells.Item(108, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(108, 5).Value = "22.9C"
$wsTemperature.Cells.Item(108, 6).Value = "Active"
$wsTemperature.Cells.Item(109, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(109, 2).Value = "16:46:09"
$wsTemperature.Cells.Item(109, 3).Value = "16:00"
$wsTemperature.Cells.Item(109, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(109, 5).Value = "22.9C"
$wsTemperature.Cells.Item(109, 6).Value = "Active"
$wsTemperature.Cells.Item(110, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(110, 2).Value = "16:46:17"
$wsTemperature.Cells.Item(110, 3).Value = "16:00"
$wsTemperature.Cells.Item(110, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(110, 5).Value = "22.9C"
$wsTemperature.Cells.Item(110, 6).Value = "Active"
$wsTemperature.Cells.Item(111, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(111, 2).Value = "16:46:21"
$wsTemperature.Cells.Item(111, 3).Value = "16:00"
$wsTemperature.Cells.Item(111, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(111, 5).Value = "22.9C"
$wsTemperature.Cells.Item(111, 6).Value = "Active"
$wsTemperature.Cells.Item(112, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(112, 2).Value = "16:46:25"
$wsTemperature.Cells.Item(112, 3).Value = "16:00"
$wsTemperature.Cells.Item(112, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(112, 5).Value = "22.9C"
$wsTemperature.Cells.Item(112, 6).Value = "Active"
$wsTemperature.Cells.Item(113, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(113, 2).Value = "16:46:29"
$wsTemperature.Cells.Item(113, 3).Value = "16:00"
$wsTemperature.Cells.Item(113, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(113, 5).Value = "22.9C"
$wsTemperature.Cells.Item(113, 6).Value = "Active"
$wsTemperature.Cells.Item(114, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(114, 2).Value = "16:46:33"
$wsTemperature.Cells.Item(114, 3).Value = "16:00"
$wsTemperature.Cells.Item(114, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(114, 5).Value = "22.9C"
$wsTemperature.Cells.Item(114, 6).Value = "Active"
$wsTemperature.Cells.Item(115, 1).Value = "'2026-01-28"
$wsTemperature.Cells.Item(115, 2).Value = "16:46:37"
$wsTemperature.Cells.Item(115, 3).Value = "16:00"
$wsTemperature.Cells.Item(115, 4).Value = "Bathroom"
$wsTemperature.Cells.Item(115, 5).Value = "22.9C"
$wsTemperature.Cells.Item(115, 6).Value = "Active"

# --- Proximity sheet: append rows 6-7 ---
$wsProximity.Cells.Item(6, 1).Value = "'2026-01-28"
$wsProximity.Cells.Item(6, 2).Value = "16:45:59"
$wsProximity.Cells.Item(6, 3).Value = "16:00"
$wsProximity.Cells.Item(6, 4).Value = "Bathroom Door"
$wsProximity.Cells.Item(6, 5).Value = "ENTER"
$wsProximity.Cells.Item(6, 6).Value = "User ENTERED Bathroom"
$wsProximity.Cells.Item(7, 1).Value = "'2026-01-28"
$wsProximity.Cells.Item(7, 2).Value = "16:46:01"
$wsProximity.Cells.Item(7, 3).Value = "16:00"
$wsProximity.Cells.Item(7, 4).Value = "Bathroom Door"
$wsProximity.Cells.Item(7, 5).Value = "EXIT"
$wsProximity.Cells.Item(7, 6).Value = "User EXITED Bathroom"
